# Applies the "multi-industry template restore" edits described in the
# commit diff: Product/Manufacturing/Production-Operations terminology is
# swapped back to the Data-Science/ML/DevOps/Training wording across all six
# worksheets. Pre-existing blank rows (r="13" on sheet 1, r="2" on most
# sheets, etc.) are NOT touched explicitly -- the engine's writer already
# omits truly-empty rows when it serialises, matching the diff's row
# removals automatically.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Resource Overview
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resource Overview")
$ws.Range("A2").Value = "PRODUCT Resource Management Plan Project"
$ws.Range("B6").Value = "Enterprise AI/ML Implementation"
$ws.Range("A18").Value = "Data Science/AI"
$ws.Range("G18").Value = "ML, Python, Statistics"
$ws.Range("A20").Value = "Data Engineering"
$ws.Range("A22").Value = "DevOps/Infrastructure"
$ws.Range("G23").Value = "Training, Communication"

# ---------------------------------------------------------------------
# Sheet 2: Detailed Staffing Plan
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Detailed Staffing Plan")
$ws.Range("A1").Value = "DETAILED STAFFING PLAN"

$ws.Range("B9").Value = "Lead Data Scientist"
$ws.Range("C9").Value = "Data Science/AI"
$ws.Range("K9").Value = "ML, Deep Learning, Python"
$ws.Range("P9").Value = "AI Lead"

$ws.Range("B10").Value = "Senior Data Scientist"
$ws.Range("C10").Value = "Data Science/AI"
$ws.Range("K10").Value = "ML, Statistics, R/Python"

$ws.Range("B11").Value = "Data Scientist"
$ws.Range("C11").Value = "Data Science/AI"
$ws.Range("K11").Value = "ML, Python, Visualization"

$ws.Range("B12").Value = "ML Engineer"
$ws.Range("C12").Value = "Data Science/AI"
$ws.Range("K12").Value = "MLOps, Python, Cloud"

$ws.Range("B13").Value = "Junior Data Scientist"
$ws.Range("C13").Value = "Data Science/AI"

$ws.Range("B18").Value = "Senior Data Engineer"
$ws.Range("C18").Value = "Data Engineering"
$ws.Range("K18").Value = "ETL, Spark, Cloud Platforms"

$ws.Range("B19").Value = "Data Engineer"
$ws.Range("C19").Value = "Data Engineering"
$ws.Range("K19").Value = "SQL, Python, Data Pipelines"

$ws.Range("B20").Value = "Cloud Data Engineer"
$ws.Range("C20").Value = "Data Engineering"

$ws.Range("B23").Value = "DevOps Engineer"
$ws.Range("C23").Value = "DevOps/Infrastructure"
$ws.Range("P23").Value = "DevOps Lead"

$ws.Range("C24").Value = "DevOps/Infrastructure"

$ws.Range("K25").Value = "Change Management, Training"

$ws.Range("B26").Value = "Training Specialist"
$ws.Range("K26").Value = "Training Design, Facilitation"

# ---------------------------------------------------------------------
# Sheet 3: Resource Timeline
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resource Timeline")
$ws.Range("B5").Value = "Lead Data Scientist"
$ws.Range("B7").Value = "Senior Data Engineer"
$ws.Range("B9").Value = "DevOps Engineer"

# ---------------------------------------------------------------------
# Sheet 4: Skills Matrix
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Skills Matrix")
$ws.Range("C3").Value = "Python/R"
$ws.Range("D3").Value = "Machine Learning"
$ws.Range("E3").Value = "Data Engineering"
$ws.Range("F3").Value = "Cloud Platforms"
$ws.Range("J3").Value = "DevOps"

$ws.Range("B5").Value = "Lead Data Scientist"
$ws.Range("B7").Value = "Senior Data Engineer"
$ws.Range("B9").Value = "DevOps Engineer"

# ---------------------------------------------------------------------
# Sheet 5: Cost Analysis
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cost Analysis")
$ws.Range("A6").Value = "Data Science/AI"
$ws.Range("A8").Value = "Data Engineering"
$ws.Range("A10").Value = "DevOps/Infrastructure"

# ---------------------------------------------------------------------
# Sheet 6: Resource Risk Assessment
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resource Risk Assessment")
$ws.Range("B5").Value = "Team lacks required ML expertise"
$ws.Range("F5").Value = "Training and external consultants"
